$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.883.17'
$ws.Range('E2').Value = '  -3.35%  '
$ws.Range('D3').Value = '3.360.11'
$ws.Range('E3').Value = '  -2.77%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'566.61"
$ws.Range('E5').Value = '  -2.05%  '
$ws.Range('D6').Value = "'148.53"
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('E9').Value = '  +0.66%  '
$ws.Range('D10').Value = "'0.122"
$ws.Range('E10').Value = '  -1.73%  '
$ws.Range('E11').Value = '  +0.94%  '
$ws.Range('D12').Value = '3.934.21'
$ws.Range('E12').Value = '  -2.80%  '
$ws.Range('E13').Value = '  +0.96%  '
$ws.Range('D14').Value = "'27.96"
$ws.Range('E14').Value = '  -2.06%  '
$ws.Range('D15').Value = '3.358.25'
$ws.Range('E15').Value = '  -2.76%  '
$ws.Range('D16').Value = "'0.0000169"
$ws.Range('E16').Value = '  -1.79%  '
$ws.Range('D17').Value = '61.001.95'
$ws.Range('E17').Value = '  -3.22%  '
$ws.Range('D18').Value = "'6.34"
$ws.Range('E18').Value = '  -1.63%  '
$ws.Range('D19').Value = "'14.19"
$ws.Range('E19').Value = '  -2.57%  '
$ws.Range('D20').Value = "'8.82"
$ws.Range('E20').Value = '  -4.08%  '
$ws.Range('D21').Value = "'374.04"
$ws.Range('E21').Value = '  -3.62%  '
$ws.Range('D22').Value = "'75.26"
$ws.Range('E22').Value = '  +0.82%  '
$ws.Range('E23').Value = '  -0.49%  '
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('D25').Value = '3.500.62'
$ws.Range('E25').Value = '  -2.60%  '
$ws.Range('D26').Value = "'0.0000109"
$ws.Range('E26').Value = '  -5.57%  '
$ws.Range('E27').Value = '  -4.36%  '
$ws.Range('D28').Value = "'0.996"
$ws.Range('E28').Value = '  -0.40%  '
$ws.Range('D29').Value = "'7.37"
$ws.Range('E29').Value = '  -3.98%  '
$ws.Range('D30').Value = "'0.999"
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('E31').Value = '  -1.80%  '
$ws.Range('D32').Value = "'7.68"
$ws.Range('E32').Value = '  -4.64%  '
$ws.Range('D33').Value = "'22.86"
$ws.Range('E33').Value = '  -2.05%  '
$ws.Range('E34').Value = '  -4.19%  '
$ws.Range('D35').Value = "'5.37"
$ws.Range('E35').Value = '  +0.37%  '
$ws.Range('D36').Value = "'168.71"
$ws.Range('E36').Value = '  -0.85%  '
$ws.Range('E37').Value = '  -5.59%  '
$ws.Range('E38').Value = '  -3.81%  '
$ws.Range('D39').Value = "'29.12"
$ws.Range('E39').Value = '  -9.45%  '
$ws.Range('D40').Value = '3.394.97'
$ws.Range('E40').Value = '  -2.78%  '
$ws.Range('D41').Value = "'0.0745"
$ws.Range('E41').Value = '  -3.90%  '
$ws.Range('D42').Value = "'42.30"
$ws.Range('E42').Value = '  -1.32%  '
$ws.Range('D43').Value = "'0.760"
$ws.Range('E43').Value = '  -4.27%  '
$ws.Range('E45').Value = '  -3.69%  '
$ws.Range('E46').Value = '  -6.19%  '
$ws.Range('D47').Value = '2.489.84'
$ws.Range('E47').Value = '  -3.45%  '
$ws.Range('E48').Value = '  -3.28%  '
$ws.Range('D49').Value = "'22.55"
$ws.Range('E49').Value = '  -0.94%  '
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('D51').Value = "'0.0261"
$ws.Range('E51').Value = '  -2.39%  '
